$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted at row 567 (Femacal de La Calera -
# Zanahoria), pushing the previously-existing rows 567..616 down to 568..617.
$ws.Rows.Item(567).EntireRow.Insert()

$ws.Range("A567").Value = 3
$ws.Range("B567").Value = "Femacal de La Calera"
$ws.Range("C567").Value = "Coquimbo"
$ws.Range("D567").Value = 45166
$ws.Range("E567").Value = 5
$ws.Range("F567").Value = 100114013
$ws.Range("G567").Value = "Zanahoria"
$ws.Range("H567").Value = "Sin especificar"
$ws.Range("I567").Value = "Primera"
$ws.Range("J567").Value = 150
$ws.Range("K567").Value = 7500
$ws.Range("L567").Value = 7500
$ws.Range("M567").Value = 7500
$ws.Range("N567").Value = "$/saco 20 kilos"
$ws.Range("O567").Value = "Provincia de Quillota"
$ws.Range("P567").Value = 375
$ws.Range("Q567").Value = 20
$ws.Range("R567").Value = "Hortaliza"
